$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 1.8.1 -> 1.8.2
$ws.Range("B3").Value = "1.8.2"

# Experimental: true -> (cleared)
$ws.Range("B7").Value = ""

# Date: 2023-10-31 -> 2025-11-18
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "2025-11-18"
$ws.Range("A8").Copy()
$ws.Range("B8").PasteSpecial(-4122)
